$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to text format so numeric-
# looking strings (e.g. "224.93", "1.01") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "37.019.73"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.008.35"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "224.93"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").Value = "0.599"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "54.77"
$ws.Range("E8").Value = "  -5.51%  "
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -4.86%  "
$ws.Range("D12").Value = "2.302.56"
$ws.Range("E12").Value = "  -3.20%  "
$ws.Range("D13").Value = "14.09"
$ws.Range("E13").Value = "  -4.87%  "
$ws.Range("D14").Value = "20.15"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").Value = "0.738"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("D16").Value = "5.10"
$ws.Range("E16").Value = "  -4.24%  "
$ws.Range("D17").Value = "2.019.46"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").Value = "36.933.79"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "6.15"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "68.57"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "0.0₃0813"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "221.85"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -6.74%  "
$ws.Range("D26").Value = "165.72"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "9.12"
$ws.Range("E27").Value = "  -8.02%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "1.35"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "18.65"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.124"
$ws.Range("E30").Value = "  -5.60%  "
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("D32").Value = "4.49"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "0.0610"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").Value = "4.40"
$ws.Range("E34").Value = "  -5.25%  "
$ws.Range("D35").Value = "2.34"
$ws.Range("E35").Value = "  -8.04%  "
$ws.Range("D36").Value = "1.85"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("D39").Value = "5.24"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "1.476.18"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "0.0215"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("D42").Value = "94.64"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").Value = "0.0916"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "16.24"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("E46").Value = "  -6.19%  "
$ws.Range("D47").Value = "1.01"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").Value = "7.15"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").Value = "2.91"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "2.190.61"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").Value = "44.27"
$ws.Range("E51").Value = "  -3.41%  "

# Restore default cell style (removes the temporary text-format styling,
# matching the workbook's original unstyled data cells).
$ws.Range("D2:E51").Style = "Normal"
